$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
}

# Row 37 and row 38 swap content (Kaspa <-> WEMIXToken) with updated values
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D37" "2.36"
$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D38" "0.105"
$ws.Range("E38").Value = "  -1.07%  "

# Update Price (D) and Volume(1h) (E) columns for remaining rows
$ws.Range("D2").Value = "42.225.48"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "2.268.52"
$ws.Range("E3").Value = "  -0.15%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue "D5" "305.59"
$ws.Range("E5").Value = "  +0.72%  "
Set-TextValue "D6" "97.52"
$ws.Range("E6").Value = "  +4.72%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +0.03%  "
Set-TextValue "D9" "0.492"
$ws.Range("E9").Value = "  +1.33%  "
Set-TextValue "D10" "35.65"
$ws.Range("E10").Value = "  +9.36%  "
Set-TextValue "D11" "0.0796"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  -0.98%  "
Set-TextValue "D13" "6.65"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "2.598.59"
$ws.Range("E14").Value = "  -0.98%  "
Set-TextValue "D15" "14.40"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "2.268.34"
$ws.Range("E16").Value = "  -0.39%  "
Set-TextValue "D17" "0.795"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "42.142.42"
$ws.Range("E18").Value = "  +0.91%  "
Set-TextValue "D19" "12.52"
$ws.Range("E19").Value = "  -2.93%  "
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").Value = "  +0.08%  "
Set-TextValue "D21" "5.98"
$ws.Range("E21").Value = "  +0.56%  "
Set-TextValue "D22" "67.73"
$ws.Range("E22").Value = "  +0.65%  "
Set-TextValue "D23" "238.00"
$ws.Range("E23").Value = "  -2.61%  "
$ws.Range("E24").Value = "  -0.08%  "
Set-TextValue "D25" "1.96"
$ws.Range("E25").Value = "  +1.39%  "
Set-TextValue "D26" "1.00"
$ws.Range("E26").Value = "  -0.02%  "
Set-TextValue "D27" "23.77"
$ws.Range("E27").Value = "  -1.21%  "
Set-TextValue "D28" "37.24"
$ws.Range("E28").Value = "  +6.19%  "
Set-TextValue "D29" "9.51"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  +1.53%  "
Set-TextValue "D31" "160.34"
$ws.Range("E31").Value = "  -0.05%  "
Set-TextValue "D32" "5.26"
$ws.Range("E32").Value = "  -0.19%  "
Set-TextValue "D33" "1.00"
$ws.Range("E33").Value = "  +0.04%  "
Set-TextValue "D34" "3.16"
$ws.Range("E34").Value = "  +4.21%  "
Set-TextValue "D35" "0.0740"
$ws.Range("E35").Value = "  -0.54%  "
Set-TextValue "D36" "17.26"
$ws.Range("E36").Value = "  +2.24%  "
Set-TextValue "D39" "1.83"
$ws.Range("E39").Value = "  +1.65%  "
Set-TextValue "D40" "0.114"
$ws.Range("E40").Value = "  -1.36%  "
Set-TextValue "D41" "4.07"
$ws.Range("E41").Value = "  +3.47%  "
Set-TextValue "D42" "2.44"
$ws.Range("E42").Value = "  +14.48%  "
$ws.Range("D43").Value = "1.985.99"
$ws.Range("E43").Value = "  -1.52%  "
Set-TextValue "D44" "0.0285"
$ws.Range("E44").Value = "  +1.07%  "
Set-TextValue "D45" "18.91"
$ws.Range("E45").Value = "  -4.78%  "
Set-TextValue "D46" "2.95"
$ws.Range("E46").Value = "  +1.30%  "
Set-TextValue "D47" "9.95"
$ws.Range("E47").Value = "  -5.09%  "
Set-TextValue "D48" "53.27"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  +0.84%  "
Set-TextValue "D50" "71.92"
$ws.Range("E50").Value = "  -1.80%  "
Set-TextValue "D51" "91.55"
$ws.Range("E51").Value = "  -0.54%  "
